# Auto-generated edit script: updates crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.402.20"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").Value = "2.304.85"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'318.15"
$ws.Range("E5").Value = "  +1.91%  "

$ws.Range("D6").Value = "'103.49"
$ws.Range("E6").Value = "  -2.10%  "

$ws.Range("D7").Value = "'0.631"
$ws.Range("E7").Value = "  +0.91%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "'0.611"
$ws.Range("E9").Value = "  +0.65%  "

$ws.Range("D10").Value = "'40.15"
$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("D11").Value = "'0.0912"
$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").Value = "'8.37"
$ws.Range("E12").Value = "  +1.20%  "

$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = "  +0.94%  "

$ws.Range("E14").Value = "  -0.16%  "

$ws.Range("E15").Value = "  -0.69%  "

$ws.Range("D16").Value = "2.651.11"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("D17").Value = "2.308.75"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").Value = "42.513.08"
$ws.Range("E18").Value = "  +1.34%  "

$ws.Range("E19").Value = "  -1.29%  "

$ws.Range("E20").Value = "  +1.41%  "

$ws.Range("D21").Value = "'73.11"
$ws.Range("E21").Value = "  -2.02%  "

$ws.Range("D22").Value = "'3.59"
$ws.Range("E22").Value = "  +2.59%  "

$ws.Range("D23").Value = "'278.59"
$ws.Range("E23").Value = "  +7.71%  "

$ws.Range("D24").Value = "'11.15"
$ws.Range("E24").Value = "  +20.41%  "

$ws.Range("D25").Value = "'2.28"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.49%  "

$ws.Range("D27").Value = "'10.87"
$ws.Range("E27").Value = "  -1.02%  "

$ws.Range("E28").Value = "  +5.84%  "

$ws.Range("D29").Value = "'22.86"
$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("D30").Value = "'36.09"
$ws.Range("E30").Value = "  +1.64%  "

$ws.Range("D31").Value = "'165.13"
$ws.Range("E31").Value = "  +1.52%  "

$ws.Range("E32").Value = "  -1.20%  "

$ws.Range("D33").Value = "'5.90"
$ws.Range("E33").Value = "  +1.23%  "

$ws.Range("E34").Value = "  +6.03%  "

$ws.Range("E35").Value = "  +2.09%  "

$ws.Range("E36").Value = "  -10.66%  "

$ws.Range("D37").Value = "'0.0373"
$ws.Range("E37").Value = "  +5.91%  "

$ws.Range("D38").Value = "'4.63"
$ws.Range("E38").Value = "  +2.65%  "

$ws.Range("D39").Value = "'3.76"
$ws.Range("E39").Value = "  +3.44%  "

$ws.Range("D40").Value = "'2.79"
$ws.Range("E40").Value = "  +2.01%  "

$ws.Range("E41").Value = "  +3.48%  "

$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "'70.15"
$ws.Range("E42").Value = "  +0.52%  "

$ws.Range("B43").Value = "BitcoinSV"
$ws.Range("C43").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D43").Value = "'96.44"
$ws.Range("E43").Value = "  -1.71%  "

$ws.Range("D44").Value = "'0.228"
$ws.Range("E44").Value = "  -0.80%  "

$ws.Range("E45").Value = "  +0.38%  "

$ws.Range("D46").Value = "'12.15"
$ws.Range("E46").Value = "  +0.55%  "

$ws.Range("D47").Value = "'81.54"
$ws.Range("E47").Value = "  +10.27%  "

$ws.Range("D48").Value = "'112.07"
$ws.Range("E48").Value = "  +0.71%  "

$ws.Range("D49").Value = "'8.95"
$ws.Range("E49").Value = "  +0.49%  "

$ws.Range("D50").Value = "'5.27"
$ws.Range("E50").Value = "  -1.73%  "

$ws.Range("D51").Value = "1.605.21"
$ws.Range("E51").Value = "  +4.16%  "

